$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells keep a text format so values remain strings (matching original inlineStr cells)
$ws.Range("B7:C17").NumberFormat = "@"
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "304.51"
$ws.Range("E2").Value = "0.93%"
$ws.Range("D3").Value = "35.58"
$ws.Range("E3").Value = "1.28%"
$ws.Range("D4").Value = "5.064"
$ws.Range("E4").Value = "0.46%"
$ws.Range("D5").Value = "0.08043"
$ws.Range("E5").Value = "0.77%"
$ws.Range("D6").Value = "1.910"
$ws.Range("E6").Value = "-0.31%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.184"
$ws.Range("E7").Value = "3.39%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "7.736"
$ws.Range("E8").Value = "-0.87%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9278"
$ws.Range("E9").Value = "0.62%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1383"
$ws.Range("E10").Value = "4.56%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1892"
$ws.Range("E11").Value = "2.41%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09109"
$ws.Range("E12").Value = "-5.06%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03628"
$ws.Range("E13").Value = "1.14%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09809"
$ws.Range("E14").Value = "-0.49%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001413"
$ws.Range("E15").Value = "1.89%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005912"
$ws.Range("E16").Value = "1.68%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.553"
$ws.Range("E17").Value = "1.42%"
$ws.Range("D18").Value = "2.949"
$ws.Range("E18").Value = "-1.37%"
$ws.Range("E19").Value = "1.60%"
$ws.Range("E20").Value = "2.39%"
$ws.Range("D21").Value = "4.885"
$ws.Range("E21").Value = "-3.45%"
$ws.Range("D22").Value = "0.2512"
$ws.Range("E22").Value = "4.65%"
$ws.Range("D23").Value = "0.04438"
$ws.Range("E23").Value = "-1.38%"
$ws.Range("E24").Value = "0.66%"
$ws.Range("D25").Value = "0.004780"
$ws.Range("E25").Value = "-0.11%"
$ws.Range("D26").Value = "0.0001562"
$ws.Range("E26").Value = "24.92%"
$ws.Range("E27").Value = "4.43%"
$ws.Range("D39").Value = "0.01953"
$ws.Range("E39").Value = "3.75%"
$ws.Range("D40").Value = "0.04877"
$ws.Range("E40").Value = "3.43%"
$ws.Range("D41").Value = "0.007644"
$ws.Range("E41").Value = "2.06%"
$ws.Range("D42").Value = "0.009260"
$ws.Range("E42").Value = "-8.63%"
$ws.Range("D43").Value = "0.1370"
$ws.Range("E43").Value = "3.36%"
$ws.Range("D44").Value = "0.002103"
$ws.Range("E44").Value = "-0.38%"
$ws.Range("D45").Value = "0.01140"
$ws.Range("E45").Value = "7.78%"
$ws.Range("D46").Value = "0.00006381"
$ws.Range("E46").Value = "2.38%"
$ws.Range("E47").Value = "0.15%"
$ws.Range("E50").Value = "0.15%"
$ws.Range("E51").Value = "0.15%"
